$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply a "0.00" (2 decimal place) number format to the statistical result
# cells in the Summary Output / ANOVA / Coefficients / Residual Output
# tables. These are four disjoint regions a user would pick with a
# Ctrl-click multi-selection before opening Format Cells > Number.
# (Applied one area at a time -- this runtime's Range.NumberFormat setter
# only honors the first area of a comma-joined multi-area reference.)
$ws.Range("B4:B7").NumberFormat = "0.00"
$ws.Range("C12:F14").NumberFormat = "0.00"
$ws.Range("B17:I18").NumberFormat = "0.00"
$ws.Range("B25:C44").NumberFormat = "0.00"

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("F26").Select()
